$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (Miniblock), shifting it and the
# following columns one place to the right.
$ws.Range("C:C").Insert()

# Header for the new column, matching the formatting of the other header cells
$ws.Range("C1").Value = "Suggested_Block"
$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").HorizontalAlignment = -4108
$ws.Range("C1").VerticalAlignment = -4160
$ws.Range("C1").Borders.LineStyle = 1

# Fill the new column with a suggested block value of 1 for every trial row
$ws.Range("C2:C21").Value = 1

# Update selection to match the saved workbook state
$ws.Range("E7").Select()
